# Update docker code and remeasure results
#
# Re-measured benchmark numbers for the "Docker" technology rows
# (Matrix size 500..5500, rows 2-12) plus a couple of UI/view tweaks that
# came along with the resave (scroll position, selection, and the width of
# the two leftmost columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [int]$Row,
        [double]$ReadTime,
        [double]$ComputationTime,
        [double]$WriteTime,
        [double]$ExecutionTime,
        [double]$MemoryUsage,
        [double]$InitializationTime,
        [double]$CpuUsage
    )

    $ws.Cells.Item($Row, 3).Value = $ReadTime            # C - Read time (s)
    $ws.Cells.Item($Row, 4).Value = $ComputationTime     # D - Computation time (s)
    $ws.Cells.Item($Row, 5).Value = $WriteTime           # E - Write time (s)
    $ws.Cells.Item($Row, 6).Value = $ExecutionTime       # F - Execution time (s)
    $ws.Cells.Item($Row, 7).Value = $MemoryUsage         # G - Memory usage (MB)
    $ws.Cells.Item($Row, 8).Value = $InitializationTime  # H - Initialization time (s)
    $ws.Cells.Item($Row, 9).Value = $CpuUsage            # I - CPU usage (%)
}

Set-RowValues  2  0.040315000000000004  0.37406200000000001  0.017547              0.43670199999999992  561.73690800000008   0.76539985170000002  98.462383000000003
Set-RowValues  3  0.14466099999999998   5.5420719999999992   0.067990000000000009  5.7593580000000006   565.30042800000012   0.77416018910000017  99.598020999999989
Set-RowValues  4  0.29814200000000002   25.746622000000002   0.14802100000000001   26.197103999999996   575.30367999999999   0.81716950510000286  99.126615999999984
Set-RowValues  5  0.51958900000000008   79.709877999999989   0.26037199999999999   80.49410499999999    589.32019000000014   0.77876618349999982  99.382960999999995
Set-RowValues  6  0.88073699999999988   164.71641100000002   0.38848999999999995   165.99040000000002   654.527692           0.89868549099998052  99.359617
Set-RowValues  7  1.2209080000000001    292.67709000000002   0.54593100000000006   294.449051           676.42490900000007   0.8933892905999905   99.145011000000011
Set-RowValues  8  1.6298109999999997    456.05521099999999   0.76395399999999991   458.45398399999993   689.03321500000015   1.2263776955000567   99.561093
Set-RowValues  9  2.2749229999999998    661.21831399999996   1.1561399999999999    664.7468922999999    719.03641500000015   1.1345237217002104   99.450023999999999
Set-RowValues 10  3.2225139999999994    1014.477743           1.7799119999999999    1019.4852539999999   753.07826999999997   1.0610198906998869   99.444389000000001
Set-RowValues 11  4.4215290000000005    1351.746848           2.4355809999999996    1358.6087310000003   791.0850549999999    1.1303745282995989   99.581677999999997
Set-RowValues 12  6.5283230000000003    1926.9737460000001    3.0917139999999996    1936.5988019999998   846.51950199999987   1.1106633383001581   99.683218000000011

# Widen column A (matrix size) and column B (technology) slightly.
$ws.Columns.Item(1).ColumnWidth = 4.333333333333333
$ws.Columns.Item(2).ColumnWidth = 6.833333333333333

# Scroll the view back to A1 (was frozen at topLeftCell="D1") and move the
# selection to I6.
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 1
$aw.ScrollRow = 1
[void]$ws.Range("I6").Select()

Write-Output "edit applied"
